# Generate Report for Handoff
#
# The b.md file has been handed off again: a new handoff .xlf was generated
# for both locales, the status flips from "Handed back: in sync with en-US"
# to "Ready for handoff", the content is no longer flagged as a duplicate,
# and an Error Detail note about the handback file being stale is recorded.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b8201a664a09ba5e0fda584516d18e680f0ccbd8/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ba9b13107f120f9e02fe19ac91517da6b56a23b2/e2e/b.md."

# ----- Overview sheet: row 3 is the b.md summary row -----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-05 08:51:13"

# ----- zh-cn sheet: row 3 is the b.md detail row -----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Ready for handoff"
# Leading apostrophe forces "False" to stay literal text instead of being
# auto-recognised as the Boolean FALSE; resetting the style afterwards drops
# the quote-prefix formatting flag so the cell keeps its original look.
$wsZh.Range("F3").Value = "'False"
$wsZh.Range("F3").Style = "Normal"
$wsZh.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-09-05 08:50:59"
$wsZh.Range("P3").Value = $errorDetail
$wsZh.Columns.Item(16).ColumnWidth = 39.1

# ----- de-de sheet: row 3 is the b.md detail row -----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("F3").Value = "'False"
$wsDe.Range("F3").Style = "Normal"
$wsDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDe.Range("H3").Value = "2016-09-05 08:51:13"
$wsDe.Range("P3").Value = $errorDetail
$wsDe.Columns.Item(16).ColumnWidth = 39.1
